$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,(@("2024-08-28", 1746.650024414062, 1322.349975585938, 1771.75, 2073.949951171875, 462.7999877929688, 8303.099914550781, 0, 226.6708544654981))
    ,(@("2024-08-29", 1695.900024414062, 1307.849975585938, 1788.400024414062, 2194.800048828125, 465, 8381.950073242188, 0.009496472342001468, 228.8234279656675))
    ,(@("2024-08-30", 1703.900024414062, 1306.050048828125, 1743.650024414062, 2183.5, 477.0499877929688, 8368.250061035156, -0.001634465975974491, 228.4494238581518))
    ,(@("2024-09-02", 1660.599975585938, 1280.449951171875, 1769.650024414062, 2220.5, 464.8500061035156, 8325.749969482422, -0.005078731065963999, 227.2891906722018))
    ,(@("2024-09-03", 1683.75, 1248.449951171875, 1810.949951171875, 2205.5, 469.6000061035156, 8357.449920654297, 0.003807458942205738, 228.1545849336934))
    ,(@("2024-09-04", 1699.5, 1288.900024414062, 1780.25, 2217.949951171875, 466.5499877929688, 8386.249938964844, 0.003446029420932761, 228.9408123458956))
    ,(@("2024-09-05", 1694.699951171875, 1272.300048828125, 1783.150024414062, 2236.14990234375, 467.75, 8389.549926757812, 0.0003934998142180441, 229.0309005130207))
    ,(@("2024-09-06", 1724.449951171875, 1234.300048828125, 1778.650024414062, 2224.39990234375, 463.3999938964844, 8351.999908447266, -0.004475808432915339, 228.0058020771063))
    ,(@("2024-09-09", 1789.300048828125, 1234.300048828125, 1817.949951171875, 2138.60009765625, 460.7999877929688, 8362.550109863281, 0.001263194627833399, 228.293817781405))
    ,(@("2024-09-10", 1799.949951171875, 1249.650024414062, 1830.099975585938, 2140.14990234375, 466.9500122070312, 8420.699890136719, 0.006953594239734628, 229.8812803576968))
    ,(@("2024-09-11", 1788.25, 1255, 1828.699951171875, 2119.85009765625, 459.2999877929688, 8369.700012207031, -0.006056489198650145, 228.4890068662385))
    ,(@("2024-09-12", 1816.650024414062, 1257.449951171875, 1852.949951171875, 2111, 459, 8415.049926757812, 0.005418344084571652, 229.7270389249819))
    ,(@("2024-09-13", 1888, 1264.349975585938, 1921.550048828125, 2101.35009765625, 451.75, 8530.500122070312, 0.01371949023681921, 232.8787767926465))
    ,(@("2024-09-16", 1934.900024414062, 1259.75, 1916, 2047.199951171875, 447.8999938964844, 8501.549957275391, -0.003393724210849177, 232.0884504496524))
    ,(@("2024-09-17", 1891.199951171875, 1257.550048828125, 1904.050048828125, 2034.699951171875, 452.75, 8445.75, -0.006563504014657771, 230.5651369733704))
    ,(@("2024-09-18", 1879.449951171875, 1232.050048828125, 1889.400024414062, 2042.550048828125, 448.8500061035156, 8390.000091552734, -0.006600942302017657, 229.0431898073524))
    ,(@("2024-09-19", 1866.650024414062, 1234.5, 1929.199951171875, 2025.699951171875, 450, 8406.049926757812, 0.00191297199403341, 229.4813430148779))
    ,(@("2024-09-20", 1930.099975585938, 1246.550048828125, 1928.400024414062, 2062.300048828125, 446.7999877929688, 8507.750061035156, 0.01209844518691423, 232.2577104647629))
    ,(@("2024-09-23", 2012.849975585938, 1229.900024414062, 1882.449951171875, 2039, 445.75, 8501.449951171875, -0.000740514215636784, 232.0857203284725))
    ,(@("2024-09-24", 2026, 1238.800048828125, 1914.400024414062, 2064.550048828125, 439.25, 8561.500122070312, 0.007063521075032611, 233.7250627052268))
    ,(@("2024-09-25", 2051.39990234375, 1252.949951171875, 1910.849975585938, 2112.050048828125, 441.5499877929688, 8651.899841308594, 0.01055886444540762, 236.1929339598256))
)

$startRow = 632
for ($k = 0; $k -lt $data.Count; $k++) {
    $row = $startRow + $k
    $vals = $data[$k]
    $dateCell = $ws.Cells.Item($row, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $vals[0]
    $dateCell.ClearFormats()
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
    $ws.Cells.Item($row, 8).Value = $vals[6]
    $ws.Cells.Item($row, 9).Value = $vals[7]
    $ws.Cells.Item($row, 10).Value = $vals[8]
}
